$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "OK / Error description" header cells E2:E3 get refined with an example
# of the actual error message shown to the user.
$ws.Range("E2").Value = 'OK / Error description: "Error! Check if server is up and running"'
$ws.Range("E3").Value = 'OK / Error description: "Error! Check if server is up and running"'

# Reflect the final selection left on the sheet after the edit.
$ws.Range("E2:E3").Select()

$wb.Save()
